# Added C Clamp to BOM
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row before row 16 (shifts old rows 16-21 down to 17-22).
$ws.Range("A16:D16").EntireRow.Insert()

# 2. Copy the formatting from row 15 (Qty/Part/Price/Link pattern) onto the
#    new row 16 so the styles match (s=2 for A-C, s=4/hyperlink-look for D).
$ws.Range("A15:D15").Copy()
$ws.Range("A16:D16").PasteSpecial(-4122)

# 3. Populate the new row with the C-Clamp line item.
#    Set D16 (the raw URL text) before B16 so the shared-string table gets
#    the URL string before the product-name string, matching the source order.
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 4).Value = "https://www.amazon.com/Grizzly-H0481-Aluminum-C-Clamps-Opening/dp/B0000DD147/ref=sr_1_31?s=power-hand-tools&ie=UTF8&qid=1476048331&sr=1-31"
$ws.Cells.Item(16, 2).Value = " Grizzly H0481 Aluminum C-Clamps, Set of 6, 1-Inch Opening "
$ws.Cells.Item(16, 3).Value = 13.29

# 4. Fix up the TOTAL row (now row 22): extend the SUM range to include the
#    newly inserted row.
$ws.Cells.Item(22, 3).Formula = "=SUM(C4:C16)"

# 5. Rebuild the hyperlinks collection. Row-insert does not shift the
#    existing Hyperlinks' ranges in this host, so clear them all and re-add
#    them in the original relationship order with the two that live below
#    the inserted row now pointing one row further down.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D18"), "https://github.com/arkorobotics/PID/blob/master/Hardware/Motor_Mounting_Plate.dxf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D19"), "https://github.com/arkorobotics/PID/blob/master/Hardware/Motor_Stick.dxf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D15"), "https://www.amazon.com/AmazonBasics-USB-2-0-Cable--Male/dp/B00NH11KIK/ref=sr_1_3?s=pc&ie=UTF8&qid=1476044754&sr=1-3&keywords=usb+b+cable") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D14"), "https://www.amazon.com/Honbay%C2%AE-120pcs-Multicolored-Female-Breadboard/dp/B017NEGTXC/ref=sr_1_1?ie=UTF8&qid=1476044697&sr=8-1&keywords=jumper+wires+0.1") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://www.pololu.com/product/1079") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://www.pololu.com/product/989") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "https://www.pololu.com/product/713") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "https://www.pololu.com/product/2191") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "https://www.pololu.com/product/1461") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D10"), "https://www.pololu.com/product/966") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D9"), "https://www.pololu.com/product/3081") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D11"), "https://www.pololu.com/product/351") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D12"), "https://www.pololu.com/product/2382") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D13"), "https://www.pololu.com/product/1957") | Out-Null

# 5b. Hyperlinks.Add() (re)stamps a brand-new "hyperlink look" style onto its
#     target cell instead of reusing the existing one (s=4), which would
#     otherwise bloat cellXfs. Re-paste the original hyperlink-cell format
#     (copied from the untouched D16) onto every cell we just touched so the
#     style indices stay exactly as they were.
$ws.Range("D16").Copy()
foreach ($addr in @("D18","D19","D15","D14","D4","D5","D6","D7","D8","D10","D9","D11","D12","D13")) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$ws.Application.CutCopyMode = $false

# 6. Scroll / selection bookkeeping to match the author's view state.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("B16").Select()
